$wb = $excel.ActiveWorkbook

# ---- Sheet "Summary" ----
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.6647940074906367
$wsSummary.Range("C2").Value = 0.625
$wsSummary.Range("D2").Value = 0.8239700374531835
$wsSummary.Range("E2").Value = 0.7108239095315024
$wsSummary.Range("F2").Value = 0.7746478873239436
$wsSummary.Range("G2").Value = 0.8140031307812723
$wsSummary.Range("H2").Value = 0.7019841770820182
$wsSummary.Range("I2").Value = 440
$wsSummary.Range("J2").Value = 264
$wsSummary.Range("K2").Value = 270
$wsSummary.Range("L2").Value = 94

# ---- Sheet "Classification Report" ----
$wsClass = $wb.Worksheets.Item("Classification Report")

# Row 2 - label "0"
$wsClass.Range("B2").Value = 0.7417582417582418
$wsClass.Range("C2").Value = 0.5056179775280899
$wsClass.Range("D2").Value = 0.6013363028953229

# Row 3 - label "1"
$wsClass.Range("B3").Value = 0.625
$wsClass.Range("C3").Value = 0.8239700374531835
$wsClass.Range("D3").Value = 0.7108239095315024

# Row 4 - label "accuracy"
$wsClass.Range("B4").Value = 0.6647940074906367
$wsClass.Range("C4").Value = 0.6647940074906367
$wsClass.Range("D4").Value = 0.6647940074906367
$wsClass.Range("E4").Value = 0.6647940074906367

# Row 5 - label "macro avg"
$wsClass.Range("B5").Value = 0.6833791208791209
$wsClass.Range("C5").Value = 0.6647940074906367
$wsClass.Range("D5").Value = 0.6560801062134127

# Row 6 - label "weighted avg"
$wsClass.Range("B6").Value = 0.6833791208791209
$wsClass.Range("C6").Value = 0.6647940074906367
$wsClass.Range("D6").Value = 0.6560801062134127

# ---- Sheet "Confusion Matrix" ----
$wsConf = $wb.Worksheets.Item("Confusion Matrix")

# Row 2 - "Actual 0"
$wsConf.Range("B2").Value = 270
$wsConf.Range("C2").Value = 264

# Row 3 - "Actual 1"
$wsConf.Range("B3").Value = 94
$wsConf.Range("C3").Value = 440
